$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage the new pretty-printed JSON text far away from row 1/2 first. Excel
# auto-expands a row's height when a value with embedded line breaks is
# typed directly into it; staging the text on an unrelated row keeps that
# side effect off of row 1, and Cut (used below) does not re-trigger it.
$staging = $ws.Range("A100")
$staging.Value2 = 'questions = [
    {
        "title": "A library you are building a system for wants to be able to retrieve book information. They tell you that some librarians remember books by their name and others by the book\u2019s numeric serial number. Which code snippet represents the best way to design endpoint(s) to handle this requirement?",
        "ques_type": 2,
        "options": [
            "@app.get(\"book/{book_id:int}\")\nasync def get_by_id():\n  pass\n\n@app.get(\"book/{book_name:str}\")\nasync def get_by_name():\n  pass\n",
            "@app.get(\"book/{book_info:str}\")\nasync def get_book():\n  pass\n",
            "@app.get(\"book/\")\nasync def get_by_id(name:str, id:int):\n  pass\n",
            "@app.get(\"book/{book_name:str}\")\nasync def get_by_name():\n  pass\n\n@app.get(\"book/{book_id:int}\")\nasync def get_by_id():\n  pass\n"
        ],
        "score": "@app.get(\"book/{book_id:int}\")\nasync def get_by_id():\n  pass\n\n@app.get(\"book/{book_name:str}\")\nasync def get_by_name():\n  pass"
    },
    {
        "title": "You are working on an inventory management system for a warehouse. The warehouse can receive new inventory for an item, but they cannot accept more than 10 boxes at a time. They need an endpoint to add the number of boxes they received. The count of boxes is passed as a query parameter.What validation would you use to minimize issues and adhere to warehouse requirements?",
        "ques_type": 2,
        "options": [
            "@app.put(\"/item/{id}\")\nasync def get_by_id(id:int,count: int=10):\n\u00a0\u00a0pass\n",
            "@app.put(\"/item/{id}\")\nasync def get_by_id(\n\u00a0\u00a0id:int,\n\u00a0\u00a0count: Annotated[int | None, Query(gt=0,le=10)] = None\n):\n\u00a0\u00a0pass\n",
            "@app.put(\"/item/{id}\")\nasync def inv(id:int, count: int = Query(...,le=10,gt=0)):\n\u00a0\u00a0pass\n",
            "@app.put(\"/item/{id}\")\nasync def inv(id:int, count: float= Path(...,le=10.0,ge=0.0)):\n\u00a0\u00a0pass\n"
        ],
        "score": "@app.put(\"/item/{id}\")\nasync def inv(id:int, count: int = Query(...,le=10,gt=0)):\n\u00a0\u00a0pass"
    },
    {
        "title": "You built a system for a health insurance firm. User info is used in many endpoints, and the firm decides they need a function that handles getting the info based on the user\u2019s API key and passes it to any endpoint that needs it. The function that performs this is shown below.Which of the below endpoints is most appropriate to utilize this function?class ProcessKey:\n    def __call__(self, api_key: Annotated[str, Header()]):\n       # Generate info based on key\n       return {\"data\":\"value\"}",
        "ques_type": 2,
        "options": [
            "@app.get(\"/\")\ndef process_form(\n  api_key: int|None=None,\n):\n  return {\"message\":\"hello\"}\n",
            "@app.get(\"/\")\ndef process_form(\n  api_key: int|None=None,\n  info: dict = Depends(ProcessKey)\n):\n  return {\"message\":\"hello\"}\n",
            "@app.get(\"/\")\ndef process_form(info: dict = Depends(ProcessKey)):\n  pass\n",
            "@app.get(\"/\")\ndef process_form(info: dict = Depends(ProcessKey())):\n  pass\n"
        ],
        "score": "@app.get(\"/\")\ndef process_form(info: dict = Depends(ProcessKey())):\n  pass"
    },
    {
        "title": "A warehouse for which you are developing a backend wants to slowly transition customers to a new warehouse. Requests to the old warehouse from certain IP addresses will be moved to the new one.What middleware should you use to best accomplish this?@app.get(\"/old_warehouse\")\ndef old_warehouse():\n    return {\"message\":\"old\"}\n\n@app.get(\"/new_warehouse\")\ndef new_warehouse():\n    return {\"message\":\"new\"}",
        "ques_type": 2,
        "options": [
            "@app.middleware(\"http\")\nasync def check(request:Request, call_next):\n  if request.url.path == ''/old_warehouse'' or request.client[0] in address_list:\n    return RedirectResponse(url=\"/new_warehouse\")\n  raise HTTPException(status_code=404, detail=\"user not in ip range\")\n",
            "@app.middleware(\"http\")\nasync def check(request:Request, call_next):\n  return await call_next(request)\n  if request.url.path == ''/old_warehouse'' and request.client[0] in address_list:\n    return RedirectResponse(url=\"/new_warehouse\")\nreturn\n",
            "@app.middleware\nasync def check(request:Request, call_next):\n  if request.url.path == ''/old_warehouse'' and request.client[0] in address_list:\n    return await call_next(request, url=\"/new_warehouse\")\n  return await call_next(request)\n",
            "@app.middleware(\"http\")\nasync def check(request:Request, call_next):\n  if request.url.path == ''/old_warehouse'' and request.client[0] in address_list:\n    return RedirectResponse(url=\"/new_warehouse\")\n  return await call_next(request)\n"
        ],
        "score": "@app.middleware(\"http\")\nasync def check(request:Request, call_next):\n  if request.url.path == ''/old_warehouse'' and request.client[0] in address_list:\n    return RedirectResponse(url=\"/new_warehouse\")\n  return await call_next(request)"
    }
]'

# Remove the old styled placeholder cell (A1 = 0, bold + bordered style)
# and the original shared-string cell (A2) so row 1 can be rebuilt cleanly.
$ws.Range("A1").ClearContents()
$ws.Range("A2").ClearContents()

# Move the staged text into A1 with default formatting.
$staging.Cut($ws.Range("A1"))

# Drop the now-empty staging row so no stray row definition remains.
$ws.Rows.Item(100).Delete()
